# Auto-generated Excel COM-interop script to apply market-price refresh diff
# Updates computed price/profit columns (H-N) across multiple worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3708937.8
$ws.Range("J17").Value = 3851516.2
$ws.Range("L17").Value = 11554548.6
$ws.Range("N17").Value = -11554884.6
$ws.Range("H53").Value = 2719.9167
$ws.Range("J53").Value = 3169.9
$ws.Range("L53").Value = 3169.9
$ws.Range("N53").Value = -4443.9
$ws.Range("H116").Value = 4580.095
$ws.Range("I116").Value = 2148
$ws.Range("J116").Value = 6791.091
$ws.Range("K116").Value = 2148
$ws.Range("L116").Value = 6791.091
$ws.Range("M116").Value = 1294
$ws.Range("N116").Value = -13675.091
$ws.Range("H129").Value = 201069.36
$ws.Range("I129").Value = 412.5
$ws.Range("J129").Value = 218517.78
$ws.Range("K129").Value = 1237.5
$ws.Range("L129").Value = 655553.34
$ws.Range("M129").Value = 3762.5
$ws.Range("N129").Value = -665553.34
$ws.Range("H137").Value = 1282.7391
$ws.Range("I137").Value = 1232.65
$ws.Range("K137").Value = 3697.95
$ws.Range("M137").Value = -1147.95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3600.6667
$ws.Range("I61").Value = 3662.4707
$ws.Range("K61").Value = 3662.4707
$ws.Range("M61").Value = -3450.4707
$ws.Range("H74").Value = 71430900
$ws.Range("I74").Value = 250000620
$ws.Range("J74").Value = 3001.4
$ws.Range("K74").Value = 250000620
$ws.Range("L74").Value = 3001.4
$ws.Range("M74").Value = -249999746
$ws.Range("N74").Value = -4749.4
$ws.Range("H77").Value = 71430900
$ws.Range("I77").Value = 250000620
$ws.Range("J77").Value = 3001.4
$ws.Range("K77").Value = 1250003100
$ws.Range("L77").Value = 15007
$ws.Range("M77").Value = -1249998732
$ws.Range("N77").Value = -23743
$ws.Range("H97").Value = 90910280
$ws.Range("I97").Value = 1120
$ws.Range("K97").Value = 1120
$ws.Range("M97").Value = -624
$ws.Range("H122").Value = 2014.7858
$ws.Range("I122").Value = 1335.9131
$ws.Range("J122").Value = 5137.6
$ws.Range("K122").Value = 4007.7393
$ws.Range("L122").Value = 15412.8
$ws.Range("M122").Value = -1557.7393
$ws.Range("N122").Value = -20312.8
$ws.Range("H136").Value = 3600.6667
$ws.Range("I136").Value = 3662.4707
$ws.Range("K136").Value = 10987.4121
$ws.Range("M136").Value = -8437.4121
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3502.4546
$ws.Range("I134").Value = 3502.4546
$ws.Range("K134").Value = 10507.3638
$ws.Range("M134").Value = -7972.363799999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17204.125
$ws.Range("I58").Value = 1362.6316
$ws.Range("K58").Value = 1362.6316
$ws.Range("M58").Value = -1159.6316
$ws.Range("H134").Value = 1368.3684
$ws.Range("I134").Value = 1281.8182
$ws.Range("J134").Value = 1487.375
$ws.Range("K134").Value = 3845.4546
$ws.Range("L134").Value = 4462.125
$ws.Range("M134").Value = -1310.4546
$ws.Range("N134").Value = -9532.125
$ws.Range("H136").Value = 17204.125
$ws.Range("I136").Value = 1362.6316
$ws.Range("K136").Value = 4087.8948
$ws.Range("M136").Value = -1537.8948
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 59.8
$ws.Range("I12").Value = 6.8
$ws.Range("J12").Value = 86.3
$ws.Range("K12").Value = 20.4
$ws.Range("L12").Value = 258.9
$ws.Range("M12").Value = 152.6
$ws.Range("N12").Value = -604.9
$ws.Range("H23").Value = 412.1111
$ws.Range("I23").Value = 34.75
$ws.Range("J23").Value = 714
$ws.Range("K23").Value = 104.25
$ws.Range("L23").Value = 2142
$ws.Range("M23").Value = 130.75
$ws.Range("N23").Value = -2612
$ws.Range("H51").Value = 2317
$ws.Range("J51").Value = 2863
$ws.Range("L51").Value = 8589
$ws.Range("N51").Value = -9509
$ws.Range("H129").Value = 269310.25
$ws.Range("I129").Value = 670
$ws.Range("J129").Value = 511086.5
$ws.Range("K129").Value = 2010
$ws.Range("L129").Value = 1533259.5
$ws.Range("M129").Value = 2990
$ws.Range("N129").Value = -1543259.5
$ws.Range("H131").Value = 717.91
$ws.Range("J131").Value = 723.55206
$ws.Range("L131").Value = 2170.65618
$ws.Range("N131").Value = -12250.65618
$ws.Range("H137").Value = 12349374
$ws.Range("J137").Value = 12824285
$ws.Range("L137").Value = 38472855
$ws.Range("N137").Value = -38483055
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1734
$ws.Range("J97").Value = 1794.5
$ws.Range("L97").Value = 1794.5
$ws.Range("N97").Value = -2786.5
$ws.Range("H122").Value = 78433256
$ws.Range("J122").Value = 166668160
$ws.Range("L122").Value = 500004480
$ws.Range("N122").Value = -500009380
$ws.Range("H132").Value = 22400.08
$ws.Range("I132").Value = 2158.4285
$ws.Range("K132").Value = 6475.2855
$ws.Range("M132").Value = -3945.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8278.916999999999
$ws.Range("I7").Value = 4294.1113
$ws.Range("J7").Value = 20233.334
$ws.Range("K7").Value = 4294.1113
$ws.Range("L7").Value = 20233.334
$ws.Range("M7").Value = -4182.1113
$ws.Range("N7").Value = -20457.334
$ws.Range("H100").Value = 2084.2307
$ws.Range("I100").Value = 1433.3334
$ws.Range("K100").Value = 1433.3334
$ws.Range("M100").Value = -892.3334
$ws.Range("H126").Value = 8278.916999999999
$ws.Range("I126").Value = 4294.1113
$ws.Range("J126").Value = 20233.334
$ws.Range("K126").Value = 12882.3339
$ws.Range("L126").Value = 60700.00199999999
$ws.Range("M126").Value = -10412.3339
$ws.Range("N126").Value = -65640.00199999999
$ws.Range("H132").Value = 2567.0667
$ws.Range("I132").Value = 2117.3333
$ws.Range("J132").Value = 4366
$ws.Range("K132").Value = 6351.999899999999
$ws.Range("L132").Value = 13098
$ws.Range("M132").Value = -3821.999899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 18999.5
$ws.Range("J63").Value = 18999.5
$ws.Range("L63").Value = 18999.5
$ws.Range("N63").Value = -20247.5
$ws.Range("H66").Value = 18999.5
$ws.Range("J66").Value = 18999.5
$ws.Range("L66").Value = 56998.5
$ws.Range("N66").Value = -63238.5
$ws.Range("H81").Value = 111112390
$ws.Range("I81").Value = 1675
$ws.Range("J81").Value = 200000960
$ws.Range("K81").Value = 3350
$ws.Range("L81").Value = 400001920
$ws.Range("M81").Value = -2289
$ws.Range("N81").Value = -400004042
$ws.Range("H84").Value = 111112390
$ws.Range("I84").Value = 1675
$ws.Range("J84").Value = 200000960
$ws.Range("K84").Value = 16750
$ws.Range("L84").Value = 2000009600
$ws.Range("M84").Value = -11446
$ws.Range("N84").Value = -2000020208
$ws.Range("H132").Value = 1908.1875
$ws.Range("I132").Value = 1211.0834
$ws.Range("K132").Value = 3633.2502
$ws.Range("M132").Value = -1103.2502

Write-Output "Applied all market-price updates"
